$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Case No"
$ws.Range("D1").Value = "PRR"
$ws.Range("E1").Value = "Percentage of cases"
$ws.Range("F1").Value = "Xsquared"
$ws.Range("G1").Value = "Significance"
$ws.Range("H1").Value = "Case No.MODERNA_BI"
$ws.Range("I1").Value = "PRR.MODERNA_BI"
$ws.Range("J1").Value = "Percentage of cases.MODERNA_BI"
$ws.Range("K1").Value = "Xsquared.MODERNA_BI"
$ws.Range("L1").Value = "Significance.MODERNA_BI"
